$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the hours value for day 7 (row 15) and its description
$ws.Range("C15").Value = 6.5
$ws.Range("D15").Value = "Finalizare 4.2 cu diagrame si nebunii"

# Recalculate so the SUM formula in D19 picks up the new value
$excel.Calculate()
